{"js": "// Add the word \"fed\" as the body text of the (until now empty) paragraph\n// that immediately follows the \"CAP\u00cdTULO II: MARCO CONCEPTUAL O REFERENCIAL\"\n// heading \u2014 part of the \"ADD marco conceptual, metodolog\u00eda seguida y\n// conclusiones\" edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"CAP\u00cdTULO II\" heading paragraph; the paragraph right after it\n// is the empty one that receives the new text.\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"CAP\u00cdTULO II: MARCO CONCEPTUAL O REFERENCIAL\") !== -1) {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex === -1 || headingIndex + 1 >= paragraphs.items.length) {\n  throw new Error('Could not locate the paragraph following \"CAP\u00cdTULO II: MARCO CONCEPTUAL O REFERENCIAL\".');\n}\n\nconst targetParagraph = paragraphs.items[headingIndex + 1];\n\n// Insert a run carrying the same Arial font formatting used throughout the\n// document, matching the run produced by the canonical edit.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p><w:r><w:rPr><w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/></w:rPr><w:t>fed</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\ntargetParagraph.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# Add the word \"fed\" as the body text of the (until now empty) paragraph\n# that immediately follows the \"CAP\u00cdTULO II: MARCO CONCEPTUAL O REFERENCIAL\"\n# heading \u2014 part of the \"ADD marco conceptual, metodolog\u00eda seguida y\n# conclusiones\" edit.\n\n$d = $word.ActiveDocument\n\n# Locate the \"CAP\u00cdTULO II\" heading paragraph (matched without accented\n# characters to stay encoding-safe); the paragraph right after it is the\n# empty one that receives the new text.\n$targetIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -match \"MARCO CONCEPTUAL\") {\n        $targetIndex = $i + 1\n        break\n    }\n    $i = $i + 1\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the paragraph following the 'CAPITULO II' heading.\"\n}\n\n$target = $d.Paragraphs($targetIndex).Range\n\n# Insert the new text at the very start of the (empty) paragraph, then give\n# it the same Arial formatting used throughout the document.\n$insertRange = $target.Duplicate\n$insertRange.Collapse(1)   # wdCollapseStart\n$insertRange.InsertBefore(\"fed\")\n$insertRange.Font.Name = \"Arial\"\n$insertRange.Font.NameAscii = \"Arial\"\n$insertRange.Font.NameFarEast = \"Arial\"\n$insertRange.Font.NameOther = \"Arial\"\n$insertRange.Font.NameBi = \"Arial\"\n"}
